$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 4967
$ws.Range("L3").Value = 5345
$ws.Range("J4").Value = 1877
$ws.Range("K4").Value = 1784
$ws.Range("L4").Value = 1309
$ws.Range("L5").Value = 316
$ws.Range("L6").Value = 4513
$ws.Range("J7").Value = 29353
$ws.Range("K7").Value = 27576
$ws.Range("L7").Value = 16450

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L2").Value = 137
$ws.Range("L4").Value = 63
$ws.Range("L6").Value = 123
$ws.Range("L7").Value = 533
$ws.Range("L9").Value = 95
$ws.Range("L11").Value = 269
$ws.Range("L19").Value = 452
$ws.Range("L22").Value = 49
$ws.Range("L23").Value = 179
$ws.Range("L29").Value = 900
$ws.Range("L33").Value = 757
$ws.Range("L34").Value = 96
$ws.Range("L36").Value = 215
$ws.Range("L37").Value = 623
$ws.Range("L42").Value = 537
$ws.Range("K43").Value = 224
$ws.Range("L48").Value = 213
$ws.Range("L50").Value = 84
$ws.Range("L51").Value = 208
$ws.Range("L52").Value = 333
$ws.Range("L57").Value = 58
$ws.Range("J63").Value = 227
$ws.Range("L63").Value = 48
$ws.Range("L67").Value = 573
$ws.Range("L70").Value = 28
$ws.Range("L72").Value = 64
$ws.Range("L73").Value = 128
$ws.Range("L76").Value = 253
$ws.Range("L77").Value = 108
$ws.Range("L78").Value = 210
$ws.Range("L83").Value = 364
$ws.Range("L84").Value = 160
$ws.Range("L85").Value = 838
$ws.Range("L87").Value = 49
$ws.Range("L88").Value = 176
$ws.Range("L90").Value = 169
$ws.Range("L95").Value = 227
$ws.Range("L96").Value = 190
$ws.Range("L99").Value = 287
$ws.Range("J101").Value = 29353
$ws.Range("K101").Value = 27576
$ws.Range("L101").Value = 16450

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("L2").Value = 112
$ws.Range("L3").Value = 146
$ws.Range("L4").Value = 15
$ws.Range("L7").Value = 364

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L2").Value = 206
$ws.Range("L3").Value = 263
$ws.Range("L5").Value = 17
$ws.Range("L7").Value = 757

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("L6").Value = 50
$ws.Range("L7").Value = 227

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L2").Value = 186
$ws.Range("L6").Value = 174
$ws.Range("L7").Value = 623

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("L2").Value = 79
$ws.Range("L4").Value = 22
$ws.Range("L6").Value = 60
$ws.Range("L7").Value = 287

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L3").Value = 219
$ws.Range("L7").Value = 573

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("L3").Value = 54
$ws.Range("L7").Value = 160

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L2").Value = 272
$ws.Range("L3").Value = 339
$ws.Range("L7").Value = 900

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("L3").Value = 54
$ws.Range("L7").Value = 213

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("L3").Value = 140
$ws.Range("L7").Value = 452

$ws = $wb.Worksheets.Item("River North")
$ws.Range("L3").Value = 47
$ws.Range("L6").Value = 115
$ws.Range("L7").Value = 253

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("L2").Value = 53
$ws.Range("L7").Value = 123

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("L3").Value = 179
$ws.Range("L6").Value = 150
$ws.Range("L7").Value = 537

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("L3").Value = 68
$ws.Range("L7").Value = 210

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("L2").Value = 45
$ws.Range("L5").Value = 4
$ws.Range("L7").Value = 179

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("L3").Value = 57
$ws.Range("L7").Value = 190

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("L3").Value = 66
$ws.Range("L7").Value = 215

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("L3").Value = 178
$ws.Range("L7").Value = 533

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("L3").Value = 26
$ws.Range("L6").Value = 32
$ws.Range("L7").Value = 96

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("L2").Value = 29
$ws.Range("L7").Value = 84

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("L2").Value = 101
$ws.Range("L6").Value = 64
$ws.Range("L7").Value = 269

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("L3").Value = 38
$ws.Range("L7").Value = 95

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("L6").Value = 31
$ws.Range("L7").Value = 128

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("L2").Value = 46
$ws.Range("L7").Value = 137

$ws = $wb.Worksheets.Item("O'Hare")
$ws.Range("L3").Value = 13
$ws.Range("L7").Value = 28

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("L2").Value = 51
$ws.Range("L7").Value = 176

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("L2").Value = 57
$ws.Range("L3").Value = 50
$ws.Range("L4").Value = 15
$ws.Range("L7").Value = 169

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("L4").Value = 30
$ws.Range("L7").Value = 208

$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Range("L4").Value = 5
$ws.Range("L7").Value = 58

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("K3").Value = 66
$ws.Range("K7").Value = 224

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L6").Value = 178
$ws.Range("L7").Value = 838

$ws = $wb.Worksheets.Item("Clearing")
$ws.Range("L3").Value = 19
$ws.Range("L7").Value = 49

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("L2").Value = 18
$ws.Range("L7").Value = 64

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("L5").Value = 7
$ws.Range("L7").Value = 108

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("L2").Value = 107
$ws.Range("L7").Value = 333

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("L2").Value = 22
$ws.Range("L7").Value = 63

$ws = $wb.Worksheets.Item("Ukrainian Village")
$ws.Range("L2").Value = 14
$ws.Range("L7").Value = 49
